$wb = $excel.ActiveWorkbook

# Overview sheet: update the "Ready for handoff" status cells for the
# f86032e9... row (row 3) to "Handed back: in sync with en-US"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: update Status + Latest Handback DateTime for row 3
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-12 00:33:10"

# de-de sheet: update Status + Latest Handback DateTime for row 3
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-12 00:33:15"
